$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-21 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("61×47=", $true, $false, $false, $false, $false, $true, 1, $false, "17×43=", 2) | Out-Null
$d.Content.Find.Execute("65×52=", $true, $false, $false, $false, $false, $true, 1, $false, "83×86=", 2) | Out-Null
$d.Content.Find.Execute("25×46=", $true, $false, $false, $false, $false, $true, 1, $false, "43×86=", 2) | Out-Null
$d.Content.Find.Execute("63×46=", $true, $false, $false, $false, $false, $true, 1, $false, "13×92=", 2) | Out-Null
$d.Content.Find.Execute("14×93=", $true, $false, $false, $false, $false, $true, 1, $false, "65×81=", 2) | Out-Null
$d.Content.Find.Execute("32×43=", $true, $false, $false, $false, $false, $true, 1, $false, "87×33=", 2) | Out-Null
$d.Content.Find.Execute("96×12=", $true, $false, $false, $false, $false, $true, 1, $false, "40×14=", 2) | Out-Null
$d.Content.Find.Execute("14×82=", $true, $false, $false, $false, $false, $true, 1, $false, "33×28=", 2) | Out-Null
$d.Content.Find.Execute("86×27=", $true, $false, $false, $false, $false, $true, 1, $false, "17×66=", 2) | Out-Null
$d.Content.Find.Execute("50×37=", $true, $false, $false, $false, $false, $true, 1, $false, "37×84=", 2) | Out-Null
$d.Content.Find.Execute("10×79=", $true, $false, $false, $false, $false, $true, 1, $false, "96×98=", 2) | Out-Null
$d.Content.Find.Execute("44×27=", $true, $false, $false, $false, $false, $true, 1, $false, "97×48=", 2) | Out-Null
$d.Content.Find.Execute("71×25=", $true, $false, $false, $false, $false, $true, 1, $false, "59×20=", 2) | Out-Null
$d.Content.Find.Execute("12×19=", $true, $false, $false, $false, $false, $true, 1, $false, "46×41=", 2) | Out-Null
$d.Content.Find.Execute("70×90=", $true, $false, $false, $false, $false, $true, 1, $false, "68×59=", 2) | Out-Null
$d.Content.Find.Execute("32×65=", $true, $false, $false, $false, $false, $true, 1, $false, "94×72=", 2) | Out-Null
$d.Content.Find.Execute("76×72=", $true, $false, $false, $false, $false, $true, 1, $false, "32×10=", 2) | Out-Null
$d.Content.Find.Execute("45×58=", $true, $false, $false, $false, $false, $true, 1, $false, "95×80=", 2) | Out-Null
$d.Content.Find.Execute("23×78=", $true, $false, $false, $false, $false, $true, 1, $false, "70×87=", 2) | Out-Null
$d.Content.Find.Execute("60×70=", $true, $false, $false, $false, $false, $true, 1, $false, "78×88=", 2) | Out-Null
$d.Content.Find.Execute("51×74=", $true, $false, $false, $false, $false, $true, 1, $false, "82×60=", 2) | Out-Null
$d.Content.Find.Execute("43×49=", $true, $false, $false, $false, $false, $true, 1, $false, "16×94=", 2) | Out-Null
$d.Content.Find.Execute("23×34=", $true, $false, $false, $false, $false, $true, 1, $false, "44×11=", 2) | Out-Null
$d.Content.Find.Execute("80×58=", $true, $false, $false, $false, $false, $true, 1, $false, "100×32=", 2) | Out-Null
$d.Content.Find.Execute("72×67=", $true, $false, $false, $false, $false, $true, 1, $false, "79×53=", 2) | Out-Null
$d.Content.Find.Execute("18×87=", $true, $false, $false, $false, $false, $true, 1, $false, "13×58=", 2) | Out-Null
$d.Content.Find.Execute("78×69=", $true, $false, $false, $false, $false, $true, 1, $false, "11×62=", 2) | Out-Null
$d.Content.Find.Execute("75×44=", $true, $false, $false, $false, $false, $true, 1, $false, "73×24=", 2) | Out-Null
$d.Content.Find.Execute("84×21=", $true, $false, $false, $false, $false, $true, 1, $false, "38×100=", 2) | Out-Null
$d.Content.Find.Execute("17×59=", $true, $false, $false, $false, $false, $true, 1, $false, "16×80=", 2) | Out-Null
$d.Content.Find.Execute("25×37=", $true, $false, $false, $false, $false, $true, 1, $false, "37×11=", 2) | Out-Null
$d.Content.Find.Execute("22×88=", $true, $false, $false, $false, $false, $true, 1, $false, "74×28=", 2) | Out-Null
$d.Content.Find.Execute("14×39=", $true, $false, $false, $false, $false, $true, 1, $false, "81×27=", 2) | Out-Null
$d.Content.Find.Execute("79×19=", $true, $false, $false, $false, $false, $true, 1, $false, "51×31=", 2) | Out-Null
$d.Content.Find.Execute("64×57=", $true, $false, $false, $false, $false, $true, 1, $false, "34×58=", 2) | Out-Null
$d.Content.Find.Execute("53×98=", $true, $false, $false, $false, $false, $true, 1, $false, "64×99=", 2) | Out-Null
$d.Content.Find.Execute("78×12=", $true, $false, $false, $false, $false, $true, 1, $false, "65×42=", 2) | Out-Null
$d.Content.Find.Execute("18×88=", $true, $false, $false, $false, $false, $true, 1, $false, "56×50=", 2) | Out-Null
$d.Content.Find.Execute("97×85=", $true, $false, $false, $false, $false, $true, 1, $false, "80×31=", 2) | Out-Null
$d.Content.Find.Execute("86×63=", $true, $false, $false, $false, $false, $true, 1, $false, "42×80=", 2) | Out-Null
$d.Content.Find.Execute("12×11=", $true, $false, $false, $false, $false, $true, 1, $false, "94×77=", 2) | Out-Null
$d.Content.Find.Execute("29×98=", $true, $false, $false, $false, $false, $true, 1, $false, "16×42=", 2) | Out-Null
$d.Content.Find.Execute("40×67=", $true, $false, $false, $false, $false, $true, 1, $false, "20×23=", 2) | Out-Null
$d.Content.Find.Execute("93×75=", $true, $false, $false, $false, $false, $true, 1, $false, "25×15=", 2) | Out-Null
$d.Content.Find.Execute("82×42=", $true, $false, $false, $false, $false, $true, 1, $false, "42×71=", 2) | Out-Null
$d.Content.Find.Execute("26×56=", $true, $false, $false, $false, $false, $true, 1, $false, "88×83=", 2) | Out-Null
$d.Content.Find.Execute("73×35=", $true, $false, $false, $false, $false, $true, 1, $false, "59×78=", 2) | Out-Null
$d.Content.Find.Execute("46×68=", $true, $false, $false, $false, $false, $true, 1, $false, "37×24=", 2) | Out-Null
$d.Content.Find.Execute("94×68=", $true, $false, $false, $false, $false, $true, 1, $false, "24×33=", 2) | Out-Null
$d.Content.Find.Execute("93×16=", $true, $false, $false, $false, $false, $true, 1, $false, "98×29=", 2) | Out-Null
$d.Content.Find.Execute("64×14=", $true, $false, $false, $false, $false, $true, 1, $false, "93×65=", 2) | Out-Null
$d.Content.Find.Execute("16×47=", $true, $false, $false, $false, $false, $true, 1, $false, "38×58=", 2) | Out-Null
$d.Content.Find.Execute("52×92=", $true, $false, $false, $false, $false, $true, 1, $false, "89×76=", 2) | Out-Null
$d.Content.Find.Execute("83×33=", $true, $false, $false, $false, $false, $true, 1, $false, "69×47=", 2) | Out-Null
$d.Content.Find.Execute("33×66=", $true, $false, $false, $false, $false, $true, 1, $false, "36×97=", 2) | Out-Null
$d.Content.Find.Execute("15×32=", $true, $false, $false, $false, $false, $true, 1, $false, "30×16=", 2) | Out-Null
$d.Content.Find.Execute("66×89=", $true, $false, $false, $false, $false, $true, 1, $false, "39×69=", 2) | Out-Null
$d.Content.Find.Execute("19×43=", $true, $false, $false, $false, $false, $true, 1, $false, "48×99=", 2) | Out-Null
$d.Content.Find.Execute("24×76=", $true, $false, $false, $false, $false, $true, 1, $false, "85×28=", 2) | Out-Null
$d.Content.Find.Execute("20×77=", $true, $false, $false, $false, $false, $true, 1, $false, "10×35=", 2) | Out-Null
$d.Content.Find.Execute("28×21=", $true, $false, $false, $false, $false, $true, 1, $false, "65×55=", 2) | Out-Null
$d.Content.Find.Execute("27×19=", $true, $false, $false, $false, $false, $true, 1, $false, "56×72=", 2) | Out-Null
$d.Content.Find.Execute("20×88=", $true, $false, $false, $false, $false, $true, 1, $false, "20×18=", 2) | Out-Null
$d.Content.Find.Execute("16×34=", $true, $false, $false, $false, $false, $true, 1, $false, "10×98=", 2) | Out-Null
$d.Content.Find.Execute("87×32=", $true, $false, $false, $false, $false, $true, 1, $false, "58×49=", 2) | Out-Null
$d.Content.Find.Execute("40×49=", $true, $false, $false, $false, $false, $true, 1, $false, "75×45=", 2) | Out-Null
$d.Content.Find.Execute("61×95=", $true, $false, $false, $false, $false, $true, 1, $false, "58×47=", 2) | Out-Null
$d.Content.Find.Execute("61×44=", $true, $false, $false, $false, $false, $true, 1, $false, "81×68=", 2) | Out-Null
$d.Content.Find.Execute("23×83=", $true, $false, $false, $false, $false, $true, 1, $false, "34×61=", 2) | Out-Null
$d.Content.Find.Execute("19×27=", $true, $false, $false, $false, $false, $true, 1, $false, "42×97=", 2) | Out-Null
$d.Content.Find.Execute("26×15=", $true, $false, $false, $false, $false, $true, 1, $false, "62×44=", 2) | Out-Null
$d.Content.Find.Execute("75×13=", $true, $false, $false, $false, $false, $true, 1, $false, "37×12=", 2) | Out-Null
$d.Content.Find.Execute("77×18=", $true, $false, $false, $false, $false, $true, 1, $false, "89×85=", 2) | Out-Null
$d.Content.Find.Execute("36×60=", $true, $false, $false, $false, $false, $true, 1, $false, "51×97=", 2) | Out-Null
$d.Content.Find.Execute("14×17=", $true, $false, $false, $false, $false, $true, 1, $false, "58×20=", 2) | Out-Null
$d.Content.Find.Execute("11×40=", $true, $false, $false, $false, $false, $true, 1, $false, "65×29=", 2) | Out-Null
$d.Content.Find.Execute("89×13=", $true, $false, $false, $false, $false, $true, 1, $false, "47×32=", 2) | Out-Null
$d.Content.Find.Execute("83×47=", $true, $false, $false, $false, $false, $true, 1, $false, "94×77=", 2) | Out-Null
$d.Content.Find.Execute("89×54=", $true, $false, $false, $false, $false, $true, 1, $false, "15×54=", 2) | Out-Null
$d.Content.Find.Execute("33×87=", $true, $false, $false, $false, $false, $true, 1, $false, "77×76=", 2) | Out-Null
$d.Content.Find.Execute("71×75=", $true, $false, $false, $false, $false, $true, 1, $false, "22×100=", 2) | Out-Null
$d.Content.Find.Execute("56×37=", $true, $false, $false, $false, $false, $true, 1, $false, "97×55=", 2) | Out-Null
$d.Content.Find.Execute("46×74=", $true, $false, $false, $false, $false, $true, 1, $false, "89×59=", 2) | Out-Null
$d.Content.Find.Execute("98×52=", $true, $false, $false, $false, $false, $true, 1, $false, "54×55=", 2) | Out-Null
$d.Content.Find.Execute("74×87=", $true, $false, $false, $false, $false, $true, 1, $false, "11×35=", 2) | Out-Null
$d.Content.Find.Execute("26×26=", $true, $false, $false, $false, $false, $true, 1, $false, "50×71=", 2) | Out-Null
$d.Content.Find.Execute("42×60=", $true, $false, $false, $false, $false, $true, 1, $false, "30×77=", 2) | Out-Null
$d.Content.Find.Execute("11×87=", $true, $false, $false, $false, $false, $true, 1, $false, "18×10=", 2) | Out-Null
$d.Content.Find.Execute("24×27=", $true, $false, $false, $false, $false, $true, 1, $false, "31×26=", 2) | Out-Null
$d.Content.Find.Execute("79×41=", $true, $false, $false, $false, $false, $true, 1, $false, "37×59=", 2) | Out-Null
$d.Content.Find.Execute("85×99=", $true, $false, $false, $false, $false, $true, 1, $false, "61×21=", 2) | Out-Null
$d.Content.Find.Execute("25×88=", $true, $false, $false, $false, $false, $true, 1, $false, "52×54=", 2) | Out-Null
$d.Content.Find.Execute("76×85=", $true, $false, $false, $false, $false, $true, 1, $false, "38×37=", 2) | Out-Null
$d.Content.Find.Execute("43×63=", $true, $false, $false, $false, $false, $true, 1, $false, "31×38=", 2) | Out-Null
$d.Content.Find.Execute("34×55=", $true, $false, $false, $false, $false, $true, 1, $false, "61×22=", 2) | Out-Null
$d.Content.Find.Execute("54×26=", $true, $false, $false, $false, $false, $true, 1, $false, "32×36=", 2) | Out-Null
$d.Content.Find.Execute("81×36=", $true, $false, $false, $false, $false, $true, 1, $false, "16×23=", 2) | Out-Null
$d.Content.Find.Execute("78×15=", $true, $false, $false, $false, $false, $true, 1, $false, "29×66=", 2) | Out-Null
$d.Content.Find.Execute("51×61=", $true, $false, $false, $false, $false, $true, 1, $false, "96×10=", 2) | Out-Null
$d.Content.Find.Execute("33×10=", $true, $false, $false, $false, $false, $true, 1, $false, "74×100=", 2) | Out-Null
